$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update A/B/C cell values for rows 2-49 ---
$ws.Range("A2").Value = "Advanced"
$ws.Range("B2").Value = "How class loader works?"
$ws.Range("C2").Value = "https://javarevisited.blogspot.com/2012/12/how-classloader-works-in-java.html"

$ws.Range("A3").Value = "Advanced"
$ws.Range("B3").Value = "How class path works in java?"
$ws.Range("C3").Value = "https://javarevisited.blogspot.com/2011/01/how-classpath-work-in-java.html"

$ws.Range("A4").Value = "Advanced"
$ws.Range("B4").Value = "Can one class be loaded by two different ClassLoader in Java?"
$ws.Range("C4").Value = "Yes"

$ws.Range("A5").Value = "Advanced"
$ws.Range("B5").Value = "How many class loader present in java?"
$ws.Range("C5").Value = 3

$ws.Range("A6").Value = "Immutability"
$ws.Range("B6").Value = "How to create a class as immutable in java which is having many fields?"
$ws.Range("C6").Value = "Need to use builder pattern."

$ws.Range("A7").Value = "Advanced"
$ws.Range("B7").Value = "In java what is the maximum number of parameter we can declare in a method?"
$ws.Range("C7").Value = 255

$ws.Range("A8").Value = "General"
$ws.Range("B8").Value = "How many ways we can create an object?"

$ws.Range("A9").Value = "Serialization&Deserialization"
$ws.Range("B9").Value = "How to avoid Deserialization process from creating another instance of Singleton class in java?"
$ws.Range("C9").Value = "Need to override readResolve() method and return the singleton object."

$ws.Range("A10").Value = "Serialization&Deserialization"
$ws.Range("B10").Value = "What is readResolve() method in java?"
$ws.Range("C10").Value = "https://www.math.uni-hamburg.de/doc/java/jdk1.4.1/docs/guide/serialization/spec/input.doc7.html"

$ws.Range("A11").Value = "Serialization&Deserialization"
$ws.Range("B11").Value = "What is the difference between Serializable and Externalizable interface?"

$ws.Range("A12").Value = "ObjectCreation"
$ws.Range("B12").Value = "What is lazy initialization?"

$ws.Range("A13").Value = "Encaptulation"
$ws.Range("B13").Value = "Can we access enum constructor from outside?"
$ws.Range("C13").Value = "https://www.java67.com/2018/07/java-enum-with-constructor-example.html?m=1"

$ws.Range("A14").Value = "General"
$ws.Range("B14").Value = "What is tight coupling and loose coupling?"
$ws.Range("C14").Value = "https://www.interviewsansar.com/2018/03/24/loose-coupling-and-tight-coupling-in-java/"

$ws.Range("A15").Value = "Multi-threading"
$ws.Range("B15").Value = "What is Runnable and Callable?"

$ws.Range("A16").Value = "General"
$ws.Range("B16").Value = "Why can't we make top level class as static?"

$ws.Range("A17").Value = "General"
$ws.Range("B17").Value = "Can we declare outer class static? Explain your answer why."

$ws.Range("A18").Value = "Collections"
$ws.Range("B18").Value = "How hash map works?"

$ws.Range("A19").Value = "Collections"
$ws.Range("B19").Value = "How ConcurrentHashMap works?"

$ws.Range("A20").Value = "Interface"
$ws.Range("B20").Value = "What is Functional Interface?"

$ws.Range("A21").Value = "Interface"
$ws.Range("B21").Value = "Can we create more than one method inside Functional Interface?"

$ws.Range("A22").Value = "Interface"
$ws.Range("B22").Value = "What is Lambda Function?"

$ws.Range("A23").Value = "Interface"
$ws.Range("B23").Value = "Can we declare final default method inside an interface? Explain your answer why."

$ws.Range("A24").Value = "Interface"
$ws.Range("B24").Value = "Can we declare final static method inside an interface? Explain your answer why."

$ws.Range("A25").Value = "Interface"
$ws.Range("B25").Value = "What is the difference in between Lambda Expression and Anonymous class?"

$ws.Range("A26").Value = "AbstractClass"
$ws.Range("B26").Value = "Can we create object for any abstract Class?"

$ws.Range("A27").Value = "Interface"
$ws.Range("B27").Value = "Can we create object for any interface?"

$ws.Range("A28").Value = "Interface"
$ws.Range("B28").Value = "What is markar interface?"

$ws.Range("A29").Value = "Serialization&Deserialization"
$ws.Range("B29").Value = "Can we serialize static field? Explain your answer why."
$ws.Range("C29").Value = "Static Variable. Static variables belong to a class and not to any individual instance. The concept of serialization is concerned with the object's current state. Only data associated with a specific instance of a class is serialized, therefore static member fields are ignored during serialization."

$ws.Range("A30").Value = "OOPs"
$ws.Range("B30").Value = "Why java does not support multiple inheritances?"

$ws.Range("A31").Value = "OOPs"
$ws.Range("B31").Value = "What is abstraction?"

$ws.Range("A32").Value = "OOPs"
$ws.Range("B32").Value = "What is encaptulation?"

$ws.Range("A33").Value = "OOPs"
$ws.Range("B33").Value = "What is the difference between abstraction and encaptulation?"

$ws.Range("A34").Value = "Multi-threading"
$ws.Range("B34").Value = "What is contex-switching in multi threading?"

$ws.Range("A35").Value = "Multi-threading"
$ws.Range("B35").Value = "How can we make sure main is the last thread to finish java program?"

$ws.Range("A36").Value = "Multi-threading"
$ws.Range("B36").Value = "How does thread communicate with each other?"

$ws.Range("A37").Value = "Multi-threading"
$ws.Range("B37").Value = "Why wait(), notify() and notifyAll() methods are in Object class?"

$ws.Range("A38").Value = "Multi-threading"
$ws.Range("B38").Value = "Why we call wait(), notify() and notifyAll() methods have to be called from synchronized method or block?"

$ws.Range("A39").Value = "Multi-threading"
$ws.Range("B39").Value = "Why sleep() and yeild() methods are static?"

$ws.Range("A40").Value = "Multi-threading"
$ws.Range("B40").Value = "Difference between interrupted() and isInterrupted() method."

$ws.Range("A41").Value = "Multi-threading"
$ws.Range("B41").Value = "How can we achieve thread safty in java?"

$ws.Range("A42").Value = "Multi-threading"
$ws.Range("B42").Value = "Which is more preferred, synchronized method or block?"

$ws.Range("A43").Value = "Multi-threading"
$ws.Range("B43").Value = "What is ThreadLocal?"
$ws.Range("C43").Value = "https://www.youtube.com/watch?v=sjMe9aecW_A"

$ws.Range("A44").Value = "Advanced"
$ws.Range("B44").Value = "How Java works?"

$ws.Range("A45").Value = "Advanced"
$ws.Range("B45").Value = "What is the difference between ClassNotFoundException and NoClassDefFoundError?"
$ws.Range("C45").Value = "https://javarevisited.blogspot.com/2011/01/how-classpath-work-in-java.html"

$ws.Range("A46").Value = "Collections"
$ws.Range("B46").Value = "What is the difference between poll() and peek() method?"

$ws.Range("A47").Value = "Multi-threading"
$ws.Range("B47").Value = "What is race condition?"

$ws.Range("A48").Value = "Multi-threading"
$ws.Range("B48").Value = "How to avoid dead lock?"
$ws.Range("C48").Value = "https://javarevisited.blogspot.com/2018/08/how-to-avoid-deadlock-in-java-threads.html"

$ws.Range("A49").Value = "Multi-threading"
$ws.Range("B49").Value = "What is Future object in java?"

# --- Add hyperlinks for Answer-link cells (text already set above; Address set here) ---
$ws.Hyperlinks.Add($ws.Range("C2"), "https://javarevisited.blogspot.com/2012/12/how-classloader-works-in-java.html", [Type]::Missing, [Type]::Missing, [Type]::Missing) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://javarevisited.blogspot.com/2011/01/how-classpath-work-in-java.html", [Type]::Missing, [Type]::Missing, [Type]::Missing) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C13"), "https://www.java67.com/2018/07/java-enum-with-constructor-example.html?m=1", [Type]::Missing, [Type]::Missing, [Type]::Missing) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C14"), "https://www.interviewsansar.com/2018/03/24/loose-coupling-and-tight-coupling-in-java/", [Type]::Missing, [Type]::Missing, [Type]::Missing) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C43"), "https://www.youtube.com/watch?v=sjMe9aecW_A", [Type]::Missing, [Type]::Missing, [Type]::Missing) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C45"), "https://javarevisited.blogspot.com/2011/01/how-classpath-work-in-java.html", [Type]::Missing, [Type]::Missing, [Type]::Missing) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C48"), "https://javarevisited.blogspot.com/2018/08/how-to-avoid-deadlock-in-java-threads.html", [Type]::Missing, [Type]::Missing, [Type]::Missing) | Out-Null

# --- Clear rows 34 extra (handled naturally by span / nothing to delete) ---

# --- Update AutoFilter range A1:C34 -> A1:C43 ---
if ($ws.AutoFilterMode) {
  $ws.AutoFilterMode = $false
}
$ws.Range("A1:C43").AutoFilter() | Out-Null

# --- Update hidden _FilterDatabase defined name ---
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
  $n = $names.Item($i)
  if ($n.Name -like "*_FilterDatabase*") {
    $n.RefersTo = "=CoreJavaTopicsOrQuestions!`$A`$1:`$C`$43"
  }
}

# --- Update frozen pane scroll position + selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B6").Select() | Out-Null

Write-Host "done"